# Fix exploration time & agent_step_time formula (#13)
# Updates columns G (Avg_Agent_Step_Time), H (Avg_Experiment_Time),
# M (Std_Agent_Step_Time), N (Std_Experiment_Time) for rows 2-13
# to reflect the corrected calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ G = 6.41424863;              H = 344.48750569;            M = 0.8515212965645228;  N = 88.38218234857277 }
    3  = @{ G = 10.97185791;             H = 984.2515611700002;       M = 1.549100725218475;   N = 275.9484917322532 }
    4  = @{ G = 1.86627926;              H = 55.36451846999999;       M = 0.3353193537182583;  N = 20.95915476158065 }
    5  = @{ G = 2.988098589999999;       H = 149.2233457;             M = 0.5331656409052948;  N = 54.40949202101353 }
    6  = @{ G = 0.49024171;              H = 8.056200349999999;       M = 0.1260675616210172;  N = 4.437641135113545 }
    7  = @{ G = 0.82102424;              H = 22.01964491;             M = 0.1790995858053783;  N = 10.44246761821167 }
    8  = @{ G = 0.21918958;              H = 2.33280048;              M = 0.05426064872468467; N = 1.248818513056049 }
    9  = @{ G = 0.3740235;               H = 6.99994068;              M = 0.09163754644614111; N = 3.785971300709705 }
    10 = @{ G = 0.11581932;              H = 0.9122120500000001;      M = 0.0270031369982331;  N = 0.440014622725552 }
    11 = @{ G = 0.20964824;              H = 3.08285206;              M = 0.0589490348301803;  N = 1.992381930183537 }
    12 = @{ G = 0.07555055999999999;     H = 0.48470289;              M = 0.01943454797368185; N = 0.2535836909713923 }
    13 = @{ G = 0.13468841;              H = 1.59421345;              M = 0.03798739386776572; N = 0.9918099429284268 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
}
